$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Name"
$ws.Range("A2").Value = "Testingchatbot"
$ws.Range("A3").Value = "Batch 03"
$ws.Range("A4").Value = "By Program Name"
$ws.Range("A5").Value = "Learning"
$ws.Range("A6").Value = "Testing"

$style = $wb.Styles.Add("MyCustomStyle")
$style.Font.Name = "Consolas"
$style.Font.Size = 12
$style.Font.Color = 16711722

$rng = $ws.Range("A2:A6")
$rng.Style = "MyCustomStyle"
$rng.EntireRow.RowHeight = 15.75

$ws.Range("A6").Select()
